$wb = $excel.ActiveWorkbook

# Rename sheets (new timestamp-based names)
$wb.Worksheets.Item(1).Name = "GNG_TO-16504778909226837"
$wb.Worksheets.Item(2).Name = "NB_TO-16504778930910718"
$wb.Worksheets.Item(3).Name = "RS_TO-16504778930969808"
$wb.Worksheets.Item(4).Name = "TOL_TO-16504778931559844"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1650477893218987"

# Sheet 1 (GNG) task order values
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1650477890893684.csv"
$ws1.Range("B3").Value = "GNG_stims-1650477890905681.csv"
$ws1.Range("B4").Value = "go_stims-16504778909076836.csv"
$ws1.Range("B5").Value = "GNG_stims-16504778909216821.csv"

# Sheet 2 (NB) task order values
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16504778930770135.csv"
$ws2.Range("B3").Value = "OB-16504778921819835.csv"
$ws2.Range("B4").Value = "ZB-match_2-16504778910236807.csv"
$ws2.Range("B5").Value = "TB-16504778928480136.csv"
$ws2.Range("B6").Value = "ZB-match_4-16504778915506806.csv"
$ws2.Range("B7").Value = "TB-16504778925270133.csv"
$ws2.Range("B8").Value = "ZB-match_5-1650477891532685.csv"
$ws2.Range("B9").Value = "OB-16504778917686841.csv"
$ws2.Range("B10").Value = "OB-1650477891849979.csv"

# Sheet 3 (RS) task order values
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# Sheet 4 (TOL) task order values
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16504778931229813.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778930989816.csv"
$ws4.Range("B4").Value = "MM_stims-1650477893139014.csv"
$ws4.Range("B5").Value = "ZM_stims-16504778931229813.csv"
$ws4.Range("B6").Value = "MM_stims-16504778931550152.csv"
$ws4.Range("B7").Value = "ZM_stims-16504778931399837.csv"

# Sheet 5 (vSAT) task order values
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16504778931590152.csv"
$ws5.Range("B3").Value = "vSAT_stims-16504778931870136.csv"
$ws5.Range("B4").Value = "vSAT_stims-16504778932030134.csv"
$ws5.Range("B5").Value = "SAT_stims-16504778931709802.csv"
